# Apply updates to the "Customer Quote" sheet:
#  - Update the Surcharge column (K) values from 1.0565 to 1 for the
#    relevant line-item rows (16, 17, 20, 23, 26, 27, 31).
#  - Move the active cell selection from A31 to F1 (to support longer
#    quotes scrolling back to the top of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

$rows = @(16, 17, 20, 23, 26, 27, 31)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = 1
}

$ws.Activate()
$ws.Range("F1").Select()
